$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.3642143333333334
$ws.Range("H2").Value = 1.092643
$ws.Range("I2").Value = 0.4800482050304226
$ws.Range("J2").Value = 0.4800482050304224
$ws.Range("M2").Value = 9.24193
$ws.Range("N2").Value = 27.72579
$ws.Range("O2").Value = 0.1468938537243544
$ws.Range("P2").Value = 0.1569651396557324
$ws.Range("Q2").Value = 3.366043373663334
$ws.Range("R2").Value = 30.29439036297
$ws.Range("S2").Value = 0.07051613081037779
$ws.Range("T2").Value = 0.07535083354408391
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.3642143333333334
$ws.Range("H3").Value = 1.092643
$ws.Range("I3").Value = 0.4800482050304226
$ws.Range("J3").Value = 0.4800482050304224
$ws.Range("O3").Value = 0.469548954544906
$ws.Range("P3").Value = 0.5017420086455576
$ws.Range("Q3").Value = 10.75962068516678
$ws.Range("R3").Value = 96.83658616650101
$ws.Range("S3").Value = 0.2254061328031936
$ws.Range("T3").Value = 0.2408603506386586
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.3642143333333334
$ws.Range("H4").Value = 1.092643
$ws.Range("I4").Value = 0.4800482050304226
$ws.Range("J4").Value = 0.4800482050304224
$ws.Range("M4").Value = 7.349831333333333
$ws.Range("N4").Value = 22.049494
$ws.Range("O4").Value = 0.1168203014713749
$ws.Range("P4").Value = 0.1248296948454213
$ws.Range("Q4").Value = 2.676913919182445
$ws.Range("R4").Value = 24.092225272642
$ws.Range("S4").Value = 0.05607937603244633
$ws.Range("T4").Value = 0.05992427094503987
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.3642143333333334
$ws.Range("H5").Value = 1.092643
$ws.Range("I5").Value = 0.4800482050304226
$ws.Range("J5").Value = 0.4800482050304224
$ws.Range("M5").Value = 12.1104985
$ws.Range("N5").Value = 24.220997
$ws.Range("O5").Value = 0.1924876941491673
$ws.Range("P5").Value = 0.1371233128688515
$ws.Range("Q5").Value = 4.410817137511833
$ws.Range("R5").Value = 26.464902825071
$ws.Range("S5").Value = 0.09240337206675273
$ws.Range("T5").Value = 0.06582580021051719
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.3642143333333334
$ws.Range("H6").Value = 1.092643
$ws.Range("I6").Value = 0.4800482050304226
$ws.Range("J6").Value = 0.4800482050304224
$ws.Range("M6").Value = 4.671440333333334
$ws.Range("N6").Value = 14.014321
$ws.Range("O6").Value = 0.07424919611019735
$ws.Range("P6").Value = 0.079339843984437
$ws.Range("Q6").Value = 1.701405526711445
$ws.Range("R6").Value = 15.312649740403
$ws.Range("S6").Value = 0.03564319331765207
$ws.Range("T6").Value = 0.03808694969212274
$ws.Range("G7").Value = 0.3944893333333333
$ws.Range("H7").Value = 1.183468
$ws.Range("I7").Value = 0.5199517949695774
$ws.Range("J7").Value = 0.5199517949695774
$ws.Range("M7").Value = 9.24193
$ws.Range("N7").Value = 27.72579
$ws.Range("O7").Value = 0.1468938537243544
$ws.Range("P7").Value = 0.1569651396557324
$ws.Range("Q7").Value = 3.645842804413333
$ws.Range("R7").Value = 32.81258523972
$ws.Range("S7").Value = 0.07637772291397664
$ws.Range("T7").Value = 0.08161430611164845
$ws.Range("G8").Value = 0.3944893333333333
$ws.Range("H8").Value = 1.183468
$ws.Range("I8").Value = 0.5199517949695774
$ws.Range("J8").Value = 0.5199517949695774
$ws.Range("O8").Value = 0.469548954544906
$ws.Range("P8").Value = 0.5017420086455576
$ws.Range("Q8").Value = 11.65400480580844
$ws.Range("R8").Value = 104.886043252276
$ws.Range("S8").Value = 0.2441428217417124
$ws.Range("T8").Value = 0.2608816580068989
$ws.Range("G9").Value = 0.3944893333333333
$ws.Range("H9").Value = 1.183468
$ws.Range("I9").Value = 0.5199517949695774
$ws.Range("J9").Value = 0.5199517949695774
$ws.Range("M9").Value = 7.349831333333333
$ws.Range("N9").Value = 22.049494
$ws.Range("O9").Value = 0.1168203014713749
$ws.Range("P9").Value = 0.1248296948454213
$ws.Range("Q9").Value = 2.899430062799111
$ws.Range("R9").Value = 26.094870565192
$ws.Range("S9").Value = 0.06074092543892853
$ws.Range("T9").Value = 0.06490542390038141
$ws.Range("G10").Value = 0.3944893333333333
$ws.Range("H10").Value = 1.183468
$ws.Range("I10").Value = 0.5199517949695774
$ws.Range("J10").Value = 0.5199517949695774
$ws.Range("M10").Value = 12.1104985
$ws.Range("N10").Value = 24.220997
$ws.Range("O10").Value = 0.1924876941491673
$ws.Range("P10").Value = 0.1371233128688515
$ws.Range("Q10").Value = 4.777462479599333
$ws.Range("R10").Value = 28.664774877596
$ws.Range("S10").Value = 0.1000843220824146
$ws.Range("T10").Value = 0.0712975126583343
$ws.Range("G11").Value = 0.3944893333333333
$ws.Range("H11").Value = 1.183468
$ws.Range("I11").Value = 0.5199517949695774
$ws.Range("J11").Value = 0.5199517949695774
$ws.Range("M11").Value = 4.671440333333334
$ws.Range("N11").Value = 14.014321
$ws.Range("O11").Value = 0.07424919611019735
$ws.Range("P11").Value = 0.079339843984437
$ws.Range("Q11").Value = 1.842833382803111
$ws.Range("R11").Value = 16.585500445228
$ws.Range("S11").Value = 0.03860600279254528
$ws.Range("T11").Value = 0.04125289429231425
